# Update countries & provincias Spain
#
# Fixes the Santa Lucia / Timor Oriental name swap (rows 204-205) and the
# Montserrat / Islas Malvinas name swap (rows 214-215), refreshes the
# "datos actualizados" timestamp, and updates the case counters that moved
# between the 05:36 and 06:53 snapshots (Peru row 5, Uzbekistan row 60,
# Tailandia row 131, Butan row 187, and the Montserrat/Islas Malvinas rows
# 214-215).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name corrections (swap mislabeled rows) ---------------------
$ws.Range("A204").Value = "Timor Oriental"
$ws.Range("A205").Value = "Santa Lucia"

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# --- Footer timestamp ------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 06:53"

# --- Peru (row 5) ------------------------------------------------------
$ws.Range("B5").Value = 5020359
$ws.Range("C5").Value = 2325
$ws.Range("D5").Value = 3942360
$ws.Range("E5").Value = 995908

# --- Uzbekistan (row 60) ------------------------------------------------
$ws.Range("B60").Value = 48565
$ws.Range("C60").Value = 136
$ws.Range("D60").Value = 45058
$ws.Range("E60").Value = 3104
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 403

# --- Tailandia (row 131) ------------------------------------------------
$ws.Range("B131").Value = 3490
$ws.Range("C131").Value = 10
$ws.Range("D131").Value = 3316
$ws.Range("E131").Value = 116

# --- Butan (row 187) ------------------------------------------------
$ws.Range("D187").Value = 175
$ws.Range("E187").Value = 71

# --- Islas Malvinas / Montserrat figures (rows 214-215) -------------------
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
